$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new row of data (row 4) mirroring the existing rows
$ws.Range("A4").Value = "김숙진"
$ws.Range("B4").Value = 19810814
$ws.Range("C4").Value = 6
$ws.Range("D4").Value = 3
$ws.Range("E4").Value = 1

# Update the active selection to match the edited workbook
$ws.Range("B8").Select()
